# Update cryptocurrency price/volume snapshot data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that can look numeric (e.g. "14.60", "0.0000268").
# Force it to Text format first so Excel keeps the exact literal formatting
# (trailing zeros, leading zeros, etc.) instead of re-parsing it as a float,
# then restore the Normal style so no stray number-format style is left behind.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '66.469.23'
$ws.Range("E2").Value = '  -0.48%  '

$ws.Range("D3").Value = '3.221.46'
$ws.Range("E3").Value = '  +0.37%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = '601.30'
$ws.Range("E5").Value = '  -0.09%  '

$ws.Range("D6").Value = '155.88'
$ws.Range("E6").Value = '  -1.51%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("D8").Value = '3.222.84'
$ws.Range("E8").Value = '  +0.36%  '

$ws.Range("D9").Value = '0.543'
$ws.Range("E9").Value = '  -1.83%  '

$ws.Range("E10").Value = '  +0.16%  '

$ws.Range("D11").Value = '5.72'
$ws.Range("E11").Value = '  -4.83%  '

$ws.Range("D12").Value = '0.500'
$ws.Range("E12").Value = '  -3.02%  '

$ws.Range("D13").Value = '0.0000268'
$ws.Range("E13").Value = '  +0.06%  '

$ws.Range("D14").Value = '38.74'
$ws.Range("E14").Value = '  -1.41%  '

$ws.Range("D15").Value = '3.754.11'
$ws.Range("E15").Value = '  +0.45%  '

$ws.Range("D16").Value = '66.520.83'
$ws.Range("E16").Value = '  -0.46%  '

$ws.Range("D17").Value = '3.243.75'
$ws.Range("E17").Value = '  +1.00%  '

$ws.Range("D18").Value = '7.25'
$ws.Range("E18").Value = '  -2.87%  '

$ws.Range("E19").Value = '  +0.92%  '

$ws.Range("D20").Value = '505.83'
$ws.Range("E20").Value = '  -2.47%  '

$ws.Range("D21").Value = '15.21'
$ws.Range("E21").Value = '  -1.43%  '

$ws.Range("D22").Value = '0.736'
$ws.Range("E22").Value = '  -0.86%  '

$ws.Range("D23").Value = '7.99'
$ws.Range("E23").Value = '  -2.51%  '

$ws.Range("D24").Value = '14.52'
$ws.Range("E24").Value = '  -3.53%  '

$ws.Range("D25").Value = '86.00'
$ws.Range("E25").Value = '  +0.75%  '

$ws.Range("D26").Value = '0.163'
$ws.Range("E26").Value = '  +79.87%  '

$ws.Range("E27").Value = '  -0.08%  '

$ws.Range("D28").Value = '2.99'
$ws.Range("E28").Value = '  -1.42%  '

$ws.Range("D29").Value = '9.00'
$ws.Range("E29").Value = '  -4.11%  '

$ws.Range("E30").Value = '  -3.60%  '

$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").Value = '6.96'
$ws.Range("E31").Value = '  -1.68%  '

$ws.Range("B32").Value = 'Stacks'
$ws.Range("C32").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D32").Value = '2.89'
$ws.Range("E32").Value = '  -6.57%  '

$ws.Range("D33").Value = '28.13'
$ws.Range("E33").Value = '  -0.67%  '

$ws.Range("E34").Value = '  +0.11%  '

$ws.Range("E35").Value = '  -6.65%  '

$ws.Range("D36").Value = '6.34'
$ws.Range("E36").Value = '  -4.07%  '

$ws.Range("D37").Value = '55.33'
$ws.Range("E37").Value = '  +0.51%  '

$ws.Range("D38").Value = '0.0₃0784'
$ws.Range("E38").Value = '  +13.18%  '

$ws.Range("D39").Value = '491.61'
$ws.Range("E39").Value = '  -6.83%  '

$ws.Range("D40").Value = '3.14'
$ws.Range("E40").Value = '  +6.47%  '

$ws.Range("D41").Value = '0.0418'
$ws.Range("E41").Value = '  -1.94%  '

$ws.Range("E42").Value = '  +0.35%  '

$ws.Range("D43").Value = '8.69'
$ws.Range("E43").Value = '  -2.77%  '

$ws.Range("D44").Value = '0.292'
$ws.Range("E44").Value = '  -4.96%  '

$ws.Range("D45").Value = '2.935.34'
$ws.Range("E45").Value = '  +1.46%  '

$ws.Range("D46").Value = '2.45'
$ws.Range("E46").Value = '  -2.07%  '

$ws.Range("D47").Value = '28.03'
$ws.Range("E47").Value = '  -2.94%  '

$ws.Range("D48").Value = '2.39'
$ws.Range("E48").Value = '  -1.10%  '

$ws.Range("E51").Value = '  -5.09%  '

# Restore the default style on the price column so only cell *content* changed
$priceRange.Style = "Normal"
